$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 625
$ws.Range("I12").Value = 600
$ws.Range("K12").Value = 600
$ws.Range("M12").Value = -430

$ws.Range("H113").Value = 12054.1
$ws.Range("I113").Value = 13171.223
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 13171.223
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = -9917.223
$ws.Range("N113").Value = -8508

$ws.Range("H132").Value = 1068.9048
$ws.Range("I132").Value = 977.35
$ws.Range("K132").Value = 2932.05
$ws.Range("M132").Value = -402.0500000000002

$ws.Range("H137").Value = 1853.5
$ws.Range("I137").Value = 1649.2
$ws.Range("J137").Value = 2057.8
$ws.Range("K137").Value = 4947.6
$ws.Range("L137").Value = 6173.400000000001
$ws.Range("M137").Value = -2397.6
$ws.Range("N137").Value = -11273.4

$ws.Range("H138").Value = 2769.1304
$ws.Range("J138").Value = 2376.0588
$ws.Range("L138").Value = 7128.176399999999
$ws.Range("N138").Value = -17408.1764

$ws.Range("H141").Value = 780123.2
$ws.Range("I141").Value = 967361.4
$ws.Range("J141").Value = 4422.143
$ws.Range("K141").Value = 2902084.2
$ws.Range("L141").Value = 13266.429
$ws.Range("M141").Value = -2896904.2
$ws.Range("N141").Value = -23626.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1789.4445
$ws.Range("I45").Value = 1406
$ws.Range("K45").Value = 1406
$ws.Range("M45").Value = -1029

$ws.Range("H74").Value = 1641.5
$ws.Range("I74").Value = 1532.3
$ws.Range("K74").Value = 1532.3
$ws.Range("M74").Value = -658.3

$ws.Range("H77").Value = 1641.5
$ws.Range("I77").Value = 1532.3
$ws.Range("K77").Value = 7661.5
$ws.Range("M77").Value = -3293.5

$ws.Range("H102").Value = 2456.2856
$ws.Range("I102").Value = 2088.9
$ws.Range("K102").Value = 2088.9
$ws.Range("M102").Value = -466.9000000000001

$ws.Range("H110").Value = 5456.5
$ws.Range("I110").Value = 900
$ws.Range("K110").Value = 900
$ws.Range("M110").Value = 1145

$ws.Range("H122").Value = 873.75
$ws.Range("I122").Value = 873.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2621.25
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -171.25
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 1430.5294
$ws.Range("I132").Value = 1058.1794
$ws.Range("K132").Value = 3174.5382
$ws.Range("M132").Value = -644.5382

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1321
$ws.Range("I20").Value = 1354.381
$ws.Range("J20").Value = 1087.3334
$ws.Range("K20").Value = 1354.381
$ws.Range("L20").Value = 1087.3334
$ws.Range("M20").Value = -1107.381
$ws.Range("N20").Value = -1581.3334

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H75").Value = 7911.5
$ws.Range("I75").Value = 6593
$ws.Range("J75").Value = 11867
$ws.Range("K75").Value = 6593
$ws.Range("L75").Value = 11867
$ws.Range("M75").Value = -5657
$ws.Range("N75").Value = -13739

$ws.Range("H78").Value = 7911.5
$ws.Range("I78").Value = 6593
$ws.Range("J78").Value = 11867
$ws.Range("K78").Value = 19779
$ws.Range("L78").Value = 35601
$ws.Range("M78").Value = -15099
$ws.Range("N78").Value = -44961

$ws.Range("H86").Value = 157254.61
$ws.Range("I86").Value = 2941
$ws.Range("J86").Value = 253700.62
$ws.Range("K86").Value = 2941
$ws.Range("L86").Value = 253700.62
$ws.Range("M86").Value = -1818
$ws.Range("N86").Value = -255946.62

$ws.Range("H89").Value = 157254.61
$ws.Range("I89").Value = 2941
$ws.Range("J89").Value = 253700.62
$ws.Range("K89").Value = 14705
$ws.Range("L89").Value = 1268503.1
$ws.Range("M89").Value = -9089
$ws.Range("N89").Value = -1279735.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 837.5
$ws.Range("I16").Value = 779
$ws.Range("K16").Value = 779
$ws.Range("M16").Value = -492

$ws.Range("H62").Value = 2763.5
$ws.Range("I62").Value = 2701.3333
$ws.Range("K62").Value = 2701.3333
$ws.Range("M62").Value = -2077.3333

$ws.Range("H65").Value = 2763.5
$ws.Range("I65").Value = 2701.3333
$ws.Range("K65").Value = 13506.6665
$ws.Range("M65").Value = -10386.6665

$ws.Range("H105").Value = 2118.3333
$ws.Range("I105").Value = 2042
$ws.Range("K105").Value = 2042
$ws.Range("M105").Value = -295

$ws.Range("H107").Value = 464.77777
$ws.Range("I107").Value = 386.3846
$ws.Range("K107").Value = 386.3846
$ws.Range("M107").Value = 1533.6154

$ws.Range("H113").Value = 837.5
$ws.Range("I113").Value = 779
$ws.Range("K113").Value = 779
$ws.Range("M113").Value = 1391

$ws.Range("H134").Value = 1855.9259
$ws.Range("I134").Value = 1574.3914
$ws.Range("K134").Value = 4723.174199999999
$ws.Range("M134").Value = -2188.174199999999

$ws.Range("H140").Value = 58856.285
$ws.Range("J140").Value = 58856.285
$ws.Range("L140").Value = 58856.285
$ws.Range("N140").Value = -69216.285

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 86.2
$ws.Range("J33").Value = 94
$ws.Range("L33").Value = 564
$ws.Range("N33").Value = -1130

$ws.Range("H69").Value = 2500
$ws.Range("I69").Value = 2500
$ws.Range("J69").Value = 2500
$ws.Range("K69").Value = 7500
$ws.Range("L69").Value = 7500
$ws.Range("M69").Value = -6689
$ws.Range("N69").Value = -9122

$ws.Range("H72").Value = 2500
$ws.Range("I72").Value = 2500
$ws.Range("J72").Value = 2500
$ws.Range("K72").Value = 22500
$ws.Range("L72").Value = 22500
$ws.Range("M72").Value = -18444
$ws.Range("N72").Value = -30612

$ws.Range("H131").Value = 13289.875
$ws.Range("I131").Value = 593.7143
$ws.Range("J131").Value = 14849.053
$ws.Range("K131").Value = 1781.1429
$ws.Range("L131").Value = 44547.159
$ws.Range("M131").Value = 3258.8571
$ws.Range("N131").Value = -54627.159

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 995
$ws.Range("I80").Value = 995
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 995
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = 3
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 995
$ws.Range("I83").Value = 995
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 4975
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = 17
$ws.Range("N83").ClearContents()

$ws.Range("H113").Value = 1499.3334
$ws.Range("I113").Value = 1498
$ws.Range("K113").Value = 1498
$ws.Range("M113").Value = 672

$ws.Range("H132").Value = 1869.5283
$ws.Range("I132").Value = 1477.725
$ws.Range("K132").Value = 4433.174999999999
$ws.Range("M132").Value = -1903.174999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1916.65
$ws.Range("I132").Value = 1709.6471
$ws.Range("J132").Value = 2069.652
$ws.Range("K132").Value = 5128.9413
$ws.Range("L132").Value = 6208.956
$ws.Range("M132").Value = -2598.9413
$ws.Range("N132").Value = -11268.956

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 56691.07
$ws.Range("I122").Value = 65939.664
$ws.Range("J122").Value = 1199.5
$ws.Range("K122").Value = 197818.992
$ws.Range("L122").Value = 3598.5
$ws.Range("M122").Value = -195368.992
$ws.Range("N122").Value = -8498.5

$ws.Range("H132").Value = 1285.8108
$ws.Range("I132").Value = 1003
$ws.Range("K132").Value = 3009
$ws.Range("M132").Value = -479

$ws.Range("H136").Value = 1944.4286
$ws.Range("I136").Value = 1645
$ws.Range("K136").Value = 4935
$ws.Range("M136").Value = -2385
